# Update countries & provincias Spain
# Refresh the COVID-19 statistics table on the "Pais" sheet with newly
# fetched data and update the "last updated" timestamp banner. The table
# is kept sorted by "Casos totales" (column B) descending, so a handful of
# countries change rank (and therefore row) as a side effect of picking up
# fresh totals; for those rows both the country name (column A) and the
# statistics (columns B:H) are written explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 22:15"

# Row -> updated cell values (country name in A when the row's occupant
# changed due to re-ranking; Casos totales/Nuevos casos/Casos activos/
# Recuperados/Casos criticos/Muertes hoy/Muertes in B:H).
$updates = @{
    4 = @{ B = 3013903; C = 30975; D = 1303535; E = 1577611; G = 188; H = 132757 }
    9 = @{ A = "España"; B = 298869; C = 414; D = 0; E = 0; G = 3; H = 28388 }
    10 = @{ A = "Chile"; B = 298557; C = 3025; D = 264371; E = 27802; G = 76; H = 6384 }
    18 = @{ B = 198037; C = 479; E = 6747; G = 4; H = 9090 }
    20 = @{ B = 168335; C = 176; D = 77308; E = 61107; G = 13; H = 29920 }
    27 = @{ B = 76222; C = 969; D = 21238; E = 51562; G = 79; H = 3422 }
    31 = @{ A = "Ecuador"; B = 62380; C = 422; D = 28872; E = 28687; G = 40; H = 4821 }
    32 = @{ A = "Irak"; B = 62275; C = 1796; D = 34741; E = 24967; G = 94; H = 2567 }
    33 = @{ A = "Belgica"; B = 62016; C = 107; D = 17091; E = 35154; H = 9771 }
    71 = @{ B = 10362; C = 342; D = 6628; E = 3697 }
    73 = @{ B = 8932; C = 2; E = 543 }
    91 = @{ B = 5241; C = 245; D = 1776; E = 3442; G = 4; H = 23 }
    92 = @{ B = 5054; C = 141; D = 1984; E = 3050; G = 4; H = 20 }
    93 = @{ B = 4948; C = 69; D = 1896; E = 2919; G = 3; H = 133 }
    95 = @{ B = 4542; C = 20; E = 416 }
    133 = @{ B = 1113; C = 8; D = 575 }
    148 = @{ B = 721; C = 1; D = 269 }
    154 = @{ B = 604; C = 10; D = 315; E = 275 }
    209 = @{ A = "Islas Malvinas" }
    210 = @{ A = "Groenlandia" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
